$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (changed) date column (C) from 2023-12-14 (45274)
# to 2023-12-15 (45275) for every existing data row (rows 2 through 27).
$ws.Range("C2:C27").Value2 = 45275

# Remove the two newest log entries (rows 28 and 29), which shrinks the
# used range from A1:Y29 down to A1:Y27.
$ws.Rows.Item(29).Delete()
$ws.Rows.Item(28).Delete()

# Row 27 no longer carries an explicit custom row height in the target
# file, so reset it back to the sheet's default height.
$ws.Rows.Item(27).AutoFit()
